$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix C31 value: 0.7 -> 1
$ws.Range("C31").Value = 1

# Add underline style marker cell at C37 (empty cell with underline font style,
# mirrors cells like B33/D11/C18 already in the sheet)
$ws.Range("C37").Font.Underline = $true

# New row 38
$ws.Range("A38").Value = "hacer configurable la ip de la printer"
$ws.Range("B38").Value = "Lucas"
$ws.Range("C38").Value = "en proceso"

# New row 39
$ws.Range("A39").Value = "primer cuota - 1 mes mas"
$ws.Range("B39").Value = "Agustina"
$ws.Range("C39").Value = "en proceso"

# New row 40
$ws.Range("A40").Value = "sacar cartel de cliente asociado con éxito"
$ws.Range("B40").Value = "Agustina"
$ws.Range("C40").Value = 1
$ws.Range("C40").NumberFormat = "0%"

# Update selection to match final cursor position
$ws.Range("E40").Select()
